# Trading update: 2026-02-17 15:18:22
# Append a new OPEN MarketMaking trade (trade #18) as row 19 on both the
# "All Trades" and "MarketMaking" sheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "MarketMaking")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item(19, 1).Value = 18

    # Date column: force text so "2026-02-17" is not auto-parsed into a
    # date serial number, then reset the style so no stray number format
    # is left applied to the cell.
    $ws.Cells.Item(19, 2).NumberFormat = "@"
    $ws.Cells.Item(19, 2).Value = "2026-02-17"
    $ws.Cells.Item(19, 2).Style = "Normal"

    $ws.Cells.Item(19, 3).Value = "15:18:18"
    $ws.Cells.Item(19, 4).Value = "MarketMaking"
    $ws.Cells.Item(19, 5).Value = "DOWN"
    $ws.Cells.Item(19, 6).Value = 0.13

    # Exit Price: trade is still OPEN, so this is an explicit empty string
    # (not a blank cell). Leading apostrophe forces text so it is stored
    # as an empty string instead of clearing the cell.
    $ws.Cells.Item(19, 7).Value = "'"
    $ws.Cells.Item(19, 7).Style = "Normal"

    $ws.Cells.Item(19, 8).Value = "OPEN"
    $ws.Cells.Item(19, 9).Value = 0
    $ws.Cells.Item(19, 10).Value = 0
    $ws.Cells.Item(19, 11).Value = 99.83200220162782
    $ws.Cells.Item(19, 12).Value = 0
    $ws.Cells.Item(19, 13).Value = 0
    $ws.Cells.Item(19, 14).Value = 0.6
    $ws.Cells.Item(19, 15).Value = "Normal spread capture: 19600 bps"

    # Exit Reason: also an explicit empty string while the trade is OPEN.
    $ws.Cells.Item(19, 16).Value = "'"
    $ws.Cells.Item(19, 16).Style = "Normal"

    $ws.Cells.Item(19, 17).Value = 0
}
